$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column F (fastqFileName) to fit the long file names, keep the other
# columns at their existing default width.
$ws.Columns.Item(6).ColumnWidth = 75.7

# The manual-status column (I) used to hold a bare numeric placeholder (34).
# It now holds a textual "[34]" marker instead, for rows 3 and 4.
$ws.Range("I3").Value = "[34]"
$ws.Range("I4").Value = "[34]"

# Update the active selection to I5.
$ws.Range("I5").Select() | Out-Null
